$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(1, 1).Value = 'テンプレートイベント'
$ws.Cells.Item(1, 2).Value = 'テンプレートイベント'
$ws.Cells.Item(2, 1).Value = 'メインフロア'
$ws.Cells.Item(2, 2).Value = 'メインフロア'
$ws.Cells.Item(3, 1).Value = '猫の集会'
$ws.Cells.Item(3, 2).Value = '猫の集会'
$ws.Cells.Item(4, 1).Value = '動画テストマップ'
$ws.Cells.Item(4, 2).Value = '動画テストマップ'
$ws.Cells.Item(5, 1).Value = '研究所'
$ws.Cells.Item(6, 1).Value = '研究棟2F'
$ws.Cells.Item(7, 1).Value = '研究所中央'
$ws.Cells.Item(8, 1).Value = 'バルコニー研究所侵入'
$ws.Cells.Item(9, 1).Value = '没'
$ws.Cells.Item(10, 1).Value = '生態実験室（イベ）'
$ws.Cells.Item(11, 1).Value = '研究棟1F'
$ws.Cells.Item(12, 1).Value = 'イベント用（研究所）'
$ws.Cells.Item(13, 1).Value = '研究棟1F（イベ）'
$ws.Cells.Item(14, 1).Value = '研究所通路'
$ws.Cells.Item(15, 1).Value = '研究棟2F（イベ）'
$ws.Cells.Item(16, 1).Value = 'トラの間'
$ws.Cells.Item(17, 1).Value = 'ブタの間'
$ws.Cells.Item(18, 1).Value = 'ヘビの間'
$ws.Cells.Item(19, 1).Value = 'カマキリの間'
$ws.Cells.Item(20, 1).Value = 'ミニマム部屋'
$ws.Cells.Item(21, 1).Value = '研究所最後の通路'
$ws.Cells.Item(22, 1).Value = 'VSベルチョコR2'
$ws.Cells.Item(23, 1).Value = '生態実験室（チョコレ戦）'
$ws.Cells.Item(24, 1).Value = 'VSベルチョコ'
$ws.Cells.Item(25, 1).Value = '生態実験室'
$ws.Cells.Item(26, 1).Value = '外エンディング用'
$ws.Cells.Item(27, 1).Value = 'ループ森'
$ws.Cells.Item(28, 1).Value = '地下階段'
$ws.Cells.Item(29, 1).Value = '日記確認部屋'
$ws.Cells.Item(30, 1).Value = '戦闘テスト'
$ws.Cells.Item(30, 2).Value = '戦闘テスト'
$ws.Cells.Item(31, 1).Value = 'イベントを作るところ'
$ws.Cells.Item(31, 2).Value = 'イベントを作るところ'
$ws.Cells.Item(32, 1).Value = 'シャル拠点'
$ws.Cells.Item(32, 2).Value = 'シャル拠点'
$ws.Cells.Item(33, 1).Value = 'テストマップ2'
$ws.Cells.Item(33, 2).Value = 'テストマップ2'
$ws.Cells.Item(34, 1).Value = '汎用イベント/テスト'
$ws.Cells.Item(34, 2).Value = '汎用イベント/テスト'
$ws.Cells.Item(35, 1).Value = 'テストマップ3'
$ws.Cells.Item(35, 2).Value = 'テストマップ3'
$ws.Cells.Item(36, 1).Value = 'アトリエ拠点'
$ws.Cells.Item(36, 2).Value = 'アトリエ拠点'
$ws.Cells.Item(37, 1).Value = 'ミスト屋敷'
$ws.Cells.Item(37, 2).Value = 'ミスト屋敷'
$ws.Cells.Item(38, 1).Value = '地下牢'
$ws.Cells.Item(38, 2).Value = '地下牢'
$ws.Cells.Item(39, 1).Value = '1Fまとめ'
$ws.Cells.Item(39, 2).Value = '1Fまとめ'
$ws.Cells.Item(40, 1).Value = 'ハンター2オープニング'
$ws.Cells.Item(40, 2).Value = 'ハンター2オープニング'
$ws.Cells.Item(41, 1).Value = '食堂'
$ws.Cells.Item(41, 2).Value = '食堂'
$ws.Cells.Item(42, 1).Value = '倉庫'
$ws.Cells.Item(42, 2).Value = '倉庫'
$ws.Cells.Item(43, 1).Value = '調理場'
$ws.Cells.Item(43, 2).Value = '調理場'
$ws.Cells.Item(44, 1).Value = '食糧庫'
$ws.Cells.Item(44, 2).Value = '食糧庫'
$ws.Cells.Item(45, 1).Value = '書庫'
$ws.Cells.Item(45, 2).Value = '書庫'
$ws.Cells.Item(46, 1).Value = '応接室'
$ws.Cells.Item(46, 2).Value = '応接室'
$ws.Cells.Item(47, 1).Value = '2Fメインフロア'
$ws.Cells.Item(47, 2).Value = '2Fメインフロア'
$ws.Cells.Item(48, 1).Value = 'バルコニー'
$ws.Cells.Item(48, 2).Value = 'バルコニー'
$ws.Cells.Item(49, 1).Value = '館主の部屋'
$ws.Cells.Item(49, 2).Value = '館主の部屋'
$ws.Cells.Item(50, 1).Value = 'シィナの部屋（元執事）'
$ws.Cells.Item(50, 2).Value = 'シィナの部屋（元執事）'
$ws.Cells.Item(51, 1).Value = 'リリーの部屋'
$ws.Cells.Item(51, 2).Value = 'リリーの部屋'
$ws.Cells.Item(52, 1).Value = 'ライムの部屋'
$ws.Cells.Item(52, 2).Value = 'ライムの部屋'
$ws.Cells.Item(53, 1).Value = '回想部屋'
$ws.Cells.Item(53, 2).Value = '回想部屋'
$ws.Cells.Item(54, 1).Value = '没テスト９－１'
$ws.Cells.Item(54, 2).Value = '没テスト９－１'
$ws.Cells.Item(55, 1).Value = 'メインフロア予備'
$ws.Cells.Item(55, 2).Value = 'メインフロア予備'
$ws.Cells.Item(56, 1).Value = 'イベント用マップ'
$ws.Cells.Item(56, 2).Value = 'イベント用マップ'
$ws.Cells.Item(57, 1).Value = 'ある日のロメリア'
$ws.Cells.Item(57, 2).Value = 'ある日のロメリア'
$ws.Cells.Item(58, 1).Value = 'お風呂イベント'
$ws.Cells.Item(58, 2).Value = 'お風呂イベント'
$ws.Cells.Item(59, 1).Value = '外'
$ws.Cells.Item(59, 2).Value = '外'
$ws.Cells.Item(60, 1).Value = 'ヤラレイベント'
$ws.Cells.Item(60, 2).Value = 'ヤラレイベント'
$ws.Cells.Item(61, 1).Value = 'ゲームオーバー地下牢'
$ws.Cells.Item(61, 2).Value = 'ゲームオーバー地下牢'
$ws.Cells.Item(62, 1).Value = '玄関壊しバトル用'
$ws.Cells.Item(62, 2).Value = '玄関壊しバトル用'
$ws.Cells.Item(63, 1).Value = 'エロイベ地下牢'
$ws.Cells.Item(63, 2).Value = 'エロイベ地下牢'
$ws.Cells.Item(64, 1).Value = 'クリア部屋'
$ws.Cells.Item(64, 2).Value = 'クリア部屋'
$ws.Cells.Item(65, 1).Value = 'お風呂前茶番'
$ws.Cells.Item(65, 2).Value = 'お風呂前茶番'
$ws.Cells.Item(66, 1).Value = '書庫イベント'
$ws.Cells.Item(66, 2).Value = '書庫イベント'
$ws.Cells.Item(67, 1).Value = 'バルコニー回想'
$ws.Cells.Item(67, 2).Value = 'バルコニー回想'
$ws.Cells.Item(68, 1).Value = '待ちぼうけ'
$ws.Cells.Item(68, 2).Value = '待ちぼうけ'
